# Add a "publisher" column to the editlist, inserted before the existing
# "journal" column (i.e. the new column becomes column B, and the old
# B..F columns shift right to C..G). Also re-home the 4 hyperlinks that
# lived in the old C3:F3 range so they land on the new D3:G3 range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Record the existing hyperlinks (column + target URL) before the column
# shift, since inserting a column does not relocate them automatically.
$links = @()
foreach ($hl in $ws.Hyperlinks) {
    $links += @{ Col = $hl.Range.Column; Row = $hl.Range.Row; Address = $hl.Address }
}

# Drop the old hyperlink relationships; they'll be rebuilt at their new
# (shifted) locations below.
$ws.Hyperlinks.Delete()

# Insert a new column at B; this shifts the old B:F columns to C:G,
# carrying cell values/styles along with them.
$ws.Columns.Item(2).Insert()

# Populate the new column B with header + values.
$ws.Range("B1").Value = "publisher"
$ws.Range("B2").Value = "sage"
$ws.Range("B3").Value = "sage"
$ws.Range("B4").Value = "sage"

# Recreate each hyperlink one column to the right of where it used to be,
# restoring the original cell text + the "Hyperlink" style that the Add
# call otherwise overrides with a freshly minted style.
foreach ($entry in $links) {
    $newCell = $ws.Cells.Item($entry.Row, $entry.Col + 1)
    $text = $newCell.Value2
    $ws.Hyperlinks.Add($newCell, $entry.Address)
    $newCell.Value = $text
    $newCell.Style = "Hyperlink"
}

# Keep the active-cell selection matching the post-edit state.
$ws.Range("B5").Select()
